$d = $word.ActiveDocument

# ---------------------------------------------------------------
# Helper: wrap a WordprocessingML body fragment in the pkg:package
# envelope that Range.InsertXML expects, then insert it, replacing
# the contents of the given range.
# ---------------------------------------------------------------
function Insert-BodyXml {
    param($range, [string]$bodyXml)

    $pkg = '<?xml version="1.0" standalone="yes"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body>' + $bodyXml + '</w:body>' +
           '</w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'

    $range.InsertXML($pkg)
}

# -----------------------------------------------------------------
# 1. Title
# -----------------------------------------------------------------
$d.Content.Find.Execute(
    "Unveiling Nature's Medicinal Symphony", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Unveiling the Enigmatic Symphony of Politics: A High School Perspective", 2) | Out-Null

# -----------------------------------------------------------------
# 2. Byline: "Isabella Peterson" -> three runs "Prof" / "." / " Julian Williamson"
# -----------------------------------------------------------------
$bylinePara = $d.Paragraphs.Item(2)
$bylineXml = '<w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:jc w:val="center"/></w:pPr>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="36"/></w:rPr><w:t>Prof</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="36"/></w:rPr><w:t>.</w:t></w:r>' +
    '<w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="36"/></w:rPr><w:t xml:space="preserve"> Julian Williamson</w:t></w:r>' +
    '</w:p>'
Insert-BodyXml $bylinePara.Range $bylineXml

# -----------------------------------------------------------------
# 3. Contact line: username + domain
# -----------------------------------------------------------------
$d.Content.Find.Execute(
    "isabella", $true, $false, $false, $false, $false,
    $true, 1, $false, "politics", 2) | Out-Null

$d.Content.Find.Execute(
    "peterson@sapiensciences", $true, $false, $false, $false, $false,
    $true, 1, $false, "simplified@schoolconnect", 2) | Out-Null

# -----------------------------------------------------------------
# 4. Big body paragraph (paragraph 5) - full rewrite
# -----------------------------------------------------------------
$bodyPara = $d.Paragraphs.Item(5)
$rPrSz24 = '<w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="24"/></w:rPr>'
$bodyXml = '<w:p>' +
    "<w:r>$rPrSz24<w:t>In the world of governance, politics stands as an enigmatic symphony, a tapestry woven with intrigue, power dynamics, and human interactions</w:t></w:r>" +
    "<w:r>$rPrSz24<w:t>.</w:t></w:r>" +
    "<w:r>$rPrSz24<w:t xml:space=`"preserve`"> At the heart of every political system lies the quest for order, stability, and the pursuit of a harmonious society</w:t></w:r>" +
    "<w:r>$rPrSz24<w:t>.</w:t></w:r>" +
    "<w:r>$rPrSz24<w:t xml:space=`"preserve`"> Yet, within this intricate symphony, there exist layers of complexity, challenges, and paradoxes that captivate the minds of students in high schools and beyond</w:t></w:r>" +
    "<w:r>$rPrSz24<w:t>.</w:t></w:r>" +
    "<w:r>$rPrSz24<w:br/></w:r>" +
    "<w:r>$rPrSz24<w:br/><w:t>Politics, in essence, is the art of negotiation, compromise, and decision-making</w:t></w:r>" +
    "<w:r>$rPrSz24<w:t>.</w:t></w:r>" +
    "<w:r>$rPrSz24<w:t xml:space=`"preserve`"> It delves into the study of how power is distributed, exercised, and contested within societies</w:t></w:r>" +
    "<w:r>$rPrSz24<w:t>.</w:t></w:r>" +
    "<w:r>$rPrSz24<w:t xml:space=`"preserve`"> Through the lens of politics, we seek to understand the motives of leaders, the influence of institutions, and the impact of policies on the lives of individuals and communities</w:t></w:r>" +
    "<w:r>$rPrSz24<w:t>.</w:t></w:r>" +
    "<w:r>$rPrSz24<w:t xml:space=`"preserve`"> It is a dynamic field where competing interests, ideologies, and values intertwine, shaping the course of nations and the lives of its citizens</w:t></w:r>" +
    "<w:r>$rPrSz24<w:t>.</w:t></w:r>" +
    "<w:r>$rPrSz24<w:br/></w:r>" +
    "<w:r>$rPrSz24<w:br/><w:t>The study of politics provides a crucial foundation for responsible citizenship</w:t></w:r>" +
    "<w:r>$rPrSz24<w:t>.</w:t></w:r>" +
    "<w:r>$rPrSz24<w:t xml:space=`"preserve`"> As future leaders, voters, and decision-makers, high school students play a pivotal role in shaping the political landscape</w:t></w:r>" +
    "<w:r>$rPrSz24<w:t>.</w:t></w:r>" +
    "<w:r>$rPrSz24<w:t xml:space=`"preserve`"> Understanding the nuances of politics equips them with the knowledge and skills necessary to navigate the complexities of governance, advocate for change, and make informed decisions that contribute to the betterment of society</w:t></w:r>" +
    "<w:r>$rPrSz24<w:t>.</w:t></w:r>" +
    '</w:p>'
Insert-BodyXml $bodyPara.Range $bodyXml

# -----------------------------------------------------------------
# 5. Summary heading: drop the lastRenderedPageBreak marker
# -----------------------------------------------------------------
$summaryHeading = $d.Paragraphs.Item(6)
$headingXml = '<w:p><w:r><w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/><w:sz w:val="28"/></w:rPr><w:t>Summary</w:t></w:r></w:p>'
Insert-BodyXml $summaryHeading.Range $headingXml

# -----------------------------------------------------------------
# 6. Summary paragraph text
# -----------------------------------------------------------------
$summaryPara = $d.Paragraphs.Item(7)
$rPrSum = '<w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri"/><w:color w:val="000000"/></w:rPr>'
$summaryXml = '<w:p>' +
    "<w:r>$rPrSum<w:t>Politics, a multifaceted and dynamic field of study, presents a symphony of intrigue, power dynamics, and human interactions</w:t></w:r>" +
    "<w:r>$rPrSum<w:t>.</w:t></w:r>" +
    "<w:r>$rPrSum<w:t xml:space=`"preserve`"> It delves into the art of negotiation, compromise, and decision-making, seeking to understand the distribution and exercise of power within societies</w:t></w:r>" +
    "<w:r>$rPrSum<w:t>.</w:t></w:r>" +
    "<w:r>$rPrSum<w:t xml:space=`"preserve`"> Politics provides a foundation for responsible citizenship, empowering high school students to navigate the complexities of governance, advocate for change, and contribute to the betterment of society</w:t></w:r>" +
    "<w:r>$rPrSum<w:t>.</w:t></w:r>" +
    '</w:p>'
Insert-BodyXml $summaryPara.Range $summaryXml

# -----------------------------------------------------------------
# 7. Append a trailing empty paragraph at the very end of the document
# -----------------------------------------------------------------
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()
